$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 23; D = 84797; E = 'love is in the air  beef fondue   sauces' }
    @{ Row = 24; D = 109439; E = 'berry  good sandwich spread' }
    @{ Row = 25; D = 42522; E = 'the man s  taco dip' }
    @{ Row = 26; D = 62368; E = 'the best  chocolate chip cheesecake ever' }
    @{ Row = 146; D = 38276; E = 'now and later  vegetarian empanadas' }
    @{ Row = 203; D = 62368; E = 'the best  chocolate chip cheesecake ever' }
    @{ Row = 204; D = 39363; E = 'the best  banana bread  or muffins' }
    @{ Row = 236; D = 41756; E = 'souper  easy sweet   sour meatballs' }
    @{ Row = 264; D = 112140; E = 'all in the kitchen  chili' }
    @{ Row = 274; D = 44045; E = 'mennonite  corn fritters' }
    @{ Row = 275; D = 87098; E = 'homemade  vegetable soup from a can' }
    @{ Row = 276; D = 112140; E = 'all in the kitchen  chili' }
    @{ Row = 304; D = 23933; E = 'chinese  candy' }
    @{ Row = 333; D = 112140; E = 'all in the kitchen  chili' }
    @{ Row = 336; D = 58224; E = 'immoral  sandwich filling  loose meat' }
    @{ Row = 355; D = 74805; E = 'never weep  whipped cream' }
    @{ Row = 356; D = 52804; E = 'jiffy  extra moist carrot cake' }
    @{ Row = 383; D = 75452; E = 'beat this  banana bread' }
    @{ Row = 384; D = 83062; E = 'spicy  banana bread' }
    @{ Row = 385; D = 39363; E = 'the best  banana bread  or muffins' }
    @{ Row = 386; D = 95926; E = 'say what   banana sandwich' }
    @{ Row = 475; D = 59952; E = 'global gourmet  taco casserole' }
    @{ Row = 476; D = 44123; E = 'george s at the cove  black bean soup' }
    @{ Row = 533; D = 38276; E = 'now and later  vegetarian empanadas' }
    @{ Row = 535; D = 67888; E = 'backyard style  barbecued ribs' }
    @{ Row = 536; D = 64045; E = 'some like it hot' }
    @{ Row = 563; D = 42570; E = 'pick me up  party chicken kabobs' }
    @{ Row = 564; D = 58224; E = 'immoral  sandwich filling  loose meat' }
    @{ Row = 665; D = 30131; E = 'momma s special  marinade' }
    @{ Row = 666; D = 93249; E = 'grilled  ranch bread' }
    @{ Row = 713; D = 83133; E = 'stove top  bbq  beef or pork ribs' }
    @{ Row = 714; D = 64302; E = 'red  macaroni salad' }
    @{ Row = 833; D = 32169; E = 'make that chicken dance  salsa pasta' }
    @{ Row = 834; D = 53402; E = 'killer  lasagna' }
    @{ Row = 835; D = 94710; E = 'italian  fries' }
    @{ Row = 836; D = 47366; E = 'forgotten  minestrone' }
    @{ Row = 923; D = 59534; E = 'twisted american chop suey' }
    @{ Row = 924; D = 41756; E = 'souper  easy sweet   sour meatballs' }
    @{ Row = 925; D = 112140; E = 'all in the kitchen  chili' }
    @{ Row = 1003; D = 52804; E = 'jiffy  extra moist carrot cake' }
    @{ Row = 1004; D = 26995; E = 'keep it going  german friendship cake' }
    @{ Row = 1006; D = 27087; E = 'get the sensation  brownies' }
    @{ Row = 1015; D = 26835; E = 'one bowl  perfect pound cake' }
    @{ Row = 1016; D = 75452; E = 'beat this  banana bread' }
    @{ Row = 1134; D = 62368; E = 'the best  chocolate chip cheesecake ever' }
    @{ Row = 1135; D = 38276; E = 'now and later  vegetarian empanadas' }
    @{ Row = 1136; D = 35653; E = 'make it your way  shortcakes' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
